# Applies the "execute condition" update to the test data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 & 3: Execute flag flips from TRUE to FALSE
$ws.Range("M2").Value = $false
$ws.Range("M3").Value = $false

# Row 4: fix email, country and add missing state; Execute stays TRUE
$ws.Range("C4").Value = "alicebrown.com"
$ws.Range("K4").Value = "India"
$ws.Range("L4").Value = "Kerala"

# Row 5: Execute flag flips from TRUE to FALSE
$ws.Range("M5").Value = $false

# Update the view so the active selection matches the edited cell
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M4").Select()
